$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate row 10 with new trial-configuration values (doric file conversion entry)
$ws.Range("A10").Value = "[180, 147.27, 114.54, 81.81, 49.09, 16.36, 0, -16.36, -49.09, -81.81, -114.54, -147.27]"
$ws.Range("B10").Value = "[2]"
$ws.Range("C10").Value = "[0.04]"
$ws.Range("D10").Value = 5
$ws.Range("E10").Value = 12
$ws.Range("F10").Value = 6
$ws.Range("H10").Value = "anaesthetized"

# Update the active selection as recorded in the workbook
$ws.Range("M11").Select()
